$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("caseType1")

$ws.Range("A6").Value = "32foobar"
$ws.Range("D6").Value = "Plain"
$ws.Range("F6").Value = 0
